# Apply updated cryptocurrency price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.180.92"
$ws.Range("E2").Value = "  -4.97%  "
$ws.Range("D3").Value = "2.236.15"
$ws.Range("E3").Value = "  -5.85%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.08"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.30"
$ws.Range("E6").Value = "  -8.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.585"
$ws.Range("E7").Value = "  -8.58%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("E9").Value = "  -8.58%  "
$ws.Range("E10").Value = "  -9.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.47"
$ws.Range("E11").Value = "  -2.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0828"
$ws.Range("E12").Value = "  -9.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.69"
$ws.Range("E13").Value = "  -10.32%  "
$ws.Range("D15").Value = "2.576.28"
$ws.Range("E15").Value = "  -5.87%  "
$ws.Range("E16").Value = "  -12.43%  "
$ws.Range("E17").Value = "  -6.78%  "
$ws.Range("D18").Value = "2.241.80"
$ws.Range("E18").Value = "  -5.37%  "
$ws.Range("D19").Value = "43.126.29"
$ws.Range("E19").Value = "  -4.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.59"
$ws.Range("E20").Value = "  -9.43%  "
$ws.Range("D21").Value = "0.0₃0966"
$ws.Range("E21").Value = "  -9.32%  "
$ws.Range("E22").Value = "  -11.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.48"
$ws.Range("E23").Value = "  -10.68%  "
$ws.Range("E24").Value = "  -11.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "237.07"
$ws.Range("E25").Value = "  -9.25%  "
$ws.Range("E26").Value = "  -8.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("E28").Value = "  +2.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.02"
$ws.Range("E29").Value = "  -10.25%  "
$ws.Range("E30").Value = "  -2.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.36"
$ws.Range("E31").Value = "  -17.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "35.40"
$ws.Range("E32").Value = "  -4.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.45"
$ws.Range("E33").Value = "  -8.66%  "
$ws.Range("E34").Value = "  -9.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "153.11"
$ws.Range("E35").Value = "  -8.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.73"
$ws.Range("E36").Value = "  -5.07%  "
$ws.Range("E37").Value = "  +8.36%  "
$ws.Range("E38").Value = "  +3.67%  "
$ws.Range("E39").Value = "  -8.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.44"
$ws.Range("E40").Value = "  -6.10%  "
$ws.Range("E41").Value = "  -12.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.66"
$ws.Range("E42").Value = "  -10.11%  "
$ws.Range("E43").Value = "  -8.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.98"
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("D46").Value = "1.783.92"
$ws.Range("E46").Value = "  -2.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.07"
$ws.Range("E47").Value = "  -12.66%  "
$ws.Range("E48").Value = "  -10.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.32"
$ws.Range("E49").Value = "  -11.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.65"
$ws.Range("E50").Value = "  -10.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.92"
$ws.Range("E51").Value = "  -16.51%  "
